# Applies the slide-4 edit:
#  - widen the teal background textbox ("文本框 1", shape id 2) so it spans
#    (almost) the full slide width
#  - delete the redundant UML picture ("图片 24", shape id 25) that used to
#    sit to the right of it
#
# PowerPoint stores shape geometry in EMU (1 pt = 12700 EMU) but the COM
# Shape.Left/Top/Width/Height properties are IEEE-754 *single*-precision
# points, and the host truncates (floors) the point->EMU conversion when it
# serialises back to OOXML. 10888661 / 12700 = 857.3748818897... does not
# round-trip exactly through a float32, so a literal with a few extra
# ten-thousandths (857.374909) is used below to land the saved <a:ext cx>
# on the exact target EMU value (10888661) instead of one EMU short.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# 1) Resize the teal text box (id 2, "文本框 1") from cx=5984239 EMU to
#    cx=10888661 EMU (height/position unchanged).
$tealBox = Get-ShapeById $s 2
$tealBox.Width = 857.374909

# 2) Delete the picture (id 25, "图片 24") that was laid over the right
#    half of the slide.
$pic = Get-ShapeById $s 25
$pic.Delete()
